$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight rows 37, 38 and 54 with a yellow fill (new fill #3 / style #4) ---
$ws.Range("A37:B38").Interior.Color = 65535
$ws.Range("A54:B54").Interior.Color = 65535

# --- Fill in the previously-empty rows 69-101 with new Q&A pairs ---
# Shared-string pool order matters: values are written in the exact order
# they appear in the target shared-string table (index 136 -> 201).
$ws.Range("A69").Value = "Nhà hàng Michelin Hà Nội – Gia có địa chỉ ở đâu, mức giá và giờ mở cửa như thế nào, và điều gì làm cho trải nghiệm ẩm thực tại nhà hàng này trở nên đặc biệt?"
$ws.Range("B69").Value = "Nhà hàng Michelin Hà Nội – Gia có địa chỉ tại 61 P. Văn Miếu, Văn Miếu, Đống Đa, Hà Nội. Mức giá tham khảo tại đây dao động khoảng 700.000 – 2.000.000 VNĐ cho set Tasting Menu.
Giờ mở cửa của nhà hàng:
Từ 18:00 – 21:00 (Thứ Ba – Thứ Năm)
Từ 11:00 – 13:30 và 18:00 – 21:00 (Thứ Năm – Chủ Nhật)
Điểm đặc biệt của Gia là phong cách phục vụ theo Tasting Menu 10 món ăn và đồ uống được sắp xếp trước, khách không chọn món theo menu truyền thống. Thực đơn được thay đổi theo mùa, mỗi 3 tháng một lần, nên khách hàng rất hiếm khi được ăn lại món cũ. Món ăn được trình bày tinh tế, tạo cảm giác mới mẻ và thú vị khi thưởng thức. Không gian nhà hàng nằm trong một ngôi nhà cổ gần 100 năm tuổi, mang đến cảm giác đẹp mắt, thân thiện và gần gũi như ở nhà."
$ws.Range("A70").Value = "Nhà hàng Bếp Prime ở Hà Nội có địa chỉ, mức giá và giờ mở cửa như thế nào? Ngoài ra, nhà hàng được mô tả có những điểm đặc sắc gì về không gian và món ăn?"
$ws.Range("B70").Value = "Nhà hàng Bếp Prime có địa chỉ tại Số 12 Hàng Điếu, Quận Hoàn Kiếm, Hà Nội. Mức giá tại đây dao động từ 350.000 - 450.000 VNĐ/người. Giờ mở cửa của nhà hàng là 10h30 - 14h00 và 17h00 - 22h30.
Bếp Prime là một trong số ít nhà hàng ở Hà Nội được Michelin đánh giá cao. Nhà hàng là sự kết hợp hài hòa giữa các yếu tố ẩm thực, cá tính đầu bếp và phong cách phục vụ chuyên nghiệp. Các món ăn được chế biến cân bằng hương vị, mang lại trải nghiệm ẩm thực khó quên. Không gian của nhà hàng ấm cúng và đậm chất Việt Nam, được trang trí với hình ảnh gánh hàng rong, tường vàng và đèn lồng truyền thống. Một món ăn đặc sắc nổi bật tại Bếp Prime là lẩu gà khiêu vũ bốc lửa, với hương vị thơm lừng mùi thảo mộc và nước lẩu ngọt thanh."
$ws.Range("A71").Value = "Duong’s Restaurant được Michelin Guide mô tả như thế nào, nhà hàng phục vụ các món gì trong menu, và thực khách có thêm lựa chọn trải nghiệm đặc biệt nào khi đến đây?"
$ws.Range("B71").Value = "Duong’s Restaurant được Michelin Guide liệt kê như một sự tôn vinh ẩm thực Việt Nam, mang đến phong cách trình bày món ăn nâng tầm và tinh tế. Nhà hàng được điều hành bởi Chef Hoàng Dương, người đã tham gia Iron Chef Vietnam và Top Chef Vietnam, và nhà hàng được đặt theo tên của ông.
Menu à la carte của Duong’s Restaurant bao gồm các món ăn đặc trưng của Hà Nội như: phở, bún chả, salad lá chuối, nem (spring rolls) và cá hấp gói lá chuối. Bên cạnh món Việt, nhà hàng cũng phục vụ một số món quốc tế như ức vịt (ăn kèm cam, truffle và cơm risotto thảo mộc). Ngoài ra, thực khách cũng có thể chọn set menu 4 món.
Một lựa chọn trải nghiệm đặc biệt tại Duong’s Restaurant là tham gia các lớp học nấu ăn chuyên nghiệp của nhà hàng nếu thực khách muốn quay lại vào một dịp khác."
$ws.Range("A72").Value = "Lụa Vạn Phúc có đặc điểm gì nổi bật theo mô tả?"
$ws.Range("B72").Value = "Lụa Vạn Phúc nổi tiếng từ xưa với nghề dệt lụa. Lụa ở đây mỏng, mềm, thoáng mát và rất nhẹ, khi mặc rất thoải mái, khác hẳn với các loại lụa thông thường."
$ws.Range("A73").Value = "Những thời điểm nào được xem là lý tưởng nhất để du lịch Hà Nội?"
$ws.Range("B73").Value = "Thời điểm lý tưởng nhất để du lịch Hà Nội là mùa xuân từ tháng 2 đến tháng 4 và mùa thu từ tháng 9 đến tháng 10. Trong những giai đoạn này, thời tiết thường mát mẻ, dễ chịu, ít mưa, rất thuận lợi cho việc tham quan và trải nghiệm các hoạt động ngoài trời."
$ws.Range("A74").Value = "thời điểm nào được xem là tốt nhất để du lịch Hà Nội?"
$ws.Range("B74").Value = "Thời điểm tốt nhất để du lịch Hà Nội là mùa thu (khoảng tháng 9 đến tháng 11) và mùa xuân (tháng 3 và tháng 4), khi thời tiết dễ chịu và nhiệt độ ôn hòa."
$ws.Range("A75").Value = "Thời điểm nào trong năm được xem là lý tưởng nhất để khám phá Hà Nội?"
$ws.Range("B75").Value = "Thời điểm lý tưởng nhất để khám phá Hà Nội là mùa xuân (tháng 2 đến tháng 3) hoặc mùa thu (tháng 8 đến tháng 10) vì đây là các mùa có thời tiết dễ chịu, mát mẻ và ổn định hơn so với mùa hè nóng bức và mùa đông lạnh."
$ws.Range("A76").Value = "Thời điểm nào trong năm là đẹp nhất để tham quan Thủ đô Hà Nội?"
$ws.Range("B76").Value = "Hà Nội có bốn mùa và mỗi mùa đều có nét đẹp riêng, vì vậy bạn có thể tham quan Hà Nội vào bất kỳ thời điểm nào trong năm.
Mùa xuân (tháng 2 – tháng 4): thời tiết lạnh, có mưa phùn nhẹ; hoa xuân nở rộ.
Mùa hè (tháng 5 – tháng 7): nắng nóng; người dân thường ra Hồ Tây hoặc vào công viên, trung tâm thương mại để giải nhiệt.
Mùa thu (tháng 8 – tháng 10): lãng mạn, thơ mộng nhất; tràn ngập hoa sữa và hương lúa mới.
Mùa đông (tháng 11 – tháng 1): rất lạnh nhưng mang vẻ đẹp trữ tình riêng.
Vì vậy, mỗi mùa đều có vẻ đẹp riêng và có thể tham quan Hà Nội vào bất kỳ thời điểm nào trong năm."
$ws.Range("A77").Value = "Những đặc sản nào được gợi ý để mua làm quà khi du lịch Hà Nội?"
$ws.Range("B77").Value = "những đặc sản được gợi ý để mua làm quà khi du lịch Hà Nội gồm:
Ô mai (Salted or Sugared Dried Fruits)
Cốm (Green Rice Flakes)
Trà sen (Lotus Tea)
Lụa Vạn Phúc (Van Phuc Silk)
Gốm Bát Tràng (Bat Trang Ceramics)
"
$ws.Range("A78").Value = "Khi du lịch Hà Nội, những món quà nào được gợi ý nên mua về?"
$ws.Range("B78").Value = "những món quà được gợi ý nên mua khi du lịch Hà Nội gồm:
Ô mai
Cốm
Trà sen
Lụa Vạn Phúc
Đồ gốm Bát Tràng"
$ws.Range("A79").Value = "Khi đi du lịch Hà Nội, những món quà nào được gợi ý mua về?"
$ws.Range("B79").Value = "ba món quà được gợi ý mua khi du lịch Hà Nội gồm:
Bánh cốm Hàng Than – có hương cốm thơm cùng nhân đậu xanh bùi bùi, là món quà quê đặc biệt.
Ô mai – đặc biệt là ô mai sấu, vị chua ngọt, được nhiều người chọn làm quà.
Lụa Vạn Phúc – lụa mềm mại, thoáng mát, màu sắc tươi sáng, hoa văn đẹp; có thể mua tại làng Vạn Phúc hoặc các gian hàng trên phố Hàng Gai và Đinh Liệt."
$ws.Range("A80").Value = "những món quà nào được gợi ý nên mua khi du lịch Hà Nội?"
$ws.Range("B80").Value = "các món quà nên mua khi du lịch Hà Nội gồm:
Bánh cốm
Trà sen
Mận Hà Nội
Chả cá Lã Vọng
Ô mai
Bánh chè lam
Sấu
Bánh chả
Bưởi diễn
Đồ gốm sứ
Đồ thủ công mỹ nghệ
Lụa
Bánh khúc
Chả rươi"
$ws.Range("A81").Value = "10 Công ty du lịch uy tín chất lượng ở Hà Nội"
$ws.Range("B81").Value = "10 Công ty du lịch uy tín chất lượng ở Hà Nội gồm:
Công ty Smile Travel
Công ty du lịch Vietravel
Công ty du lịch VietFlame Tour
Công ty Hoàng Việt Travel
Công ty Bạch Dương Travel
Công ty du lịch Best Price
Công ty du lịch Quốc tế Đại Việt
Công ty du lịch Khát Vọng Việt
Công ty Elite Tour
Công ty Bốn Mùa Tourist"
$ws.Range("A82").Value = "Phở Hà Nội gồm những đặc điểm gì?"
$ws.Range("B82").Value = "Phở Hà Nội có nước dùng trong, ngọt thanh, bánh phở mềm dai và thịt bò tươi ngon. Phở Hà Nội có nhiều loại như phở bò tái, phở bò chín, phở gà và phở chay."
$ws.Range("A83").Value = "Kem Tràng Tiền có những vị nào?"
$ws.Range("B83").Value = "Kem Tràng Tiền có nhiều vị như vani, sô cô la và dâu tây. Kem có vị ngọt thanh, mát lạnh và thơm ngon."
$ws.Range("A84").Value = "Địa chỉ của Quán Phở 10 Lý Quốc Sư là gì?"
$ws.Range("B84").Value = "Quán Phở 10 Lý Quốc Sư có địa chỉ: 10 Lý Quốc Sư, Phường Hàng Trống, Quận Hoàn Kiếm, Hà Nội, Việt Nam."
$ws.Range("A85").Value = "Cốm làng Vòng được làm từ loại lúa gì?"
$ws.Range("B85").Value = "Cốm làng Vòng được làm từ lúa nếp cái hoa vàng, có màu xanh non, dẻo thơm và vị ngọt thanh."
$ws.Range("A86").Value = "Quán Bún Chả Hương Liên nổi tiếng với món gì?"
$ws.Range("B86").Value = "Quán Bún Chả Hương Liên nổi tiếng với món Bún Chả “Obama” – từng được cựu Tổng thống Mỹ Barack Obama thưởng thức. Hương vị bún chả ở đây hài hòa, đậm đà với thịt nướng thơm lừng, chả viên mềm và nước chấm chua ngọt."
$ws.Range("A87").Value = "Hãy mô tả đầy đủ các thành phần và hương vị của món Bún Thang Hà Nội"
$ws.Range("B87").Value = "Bún thang gồm bún, giò lụa, thịt gà, trứng gà, tôm khô và các nguyên liệu khác. Nước dùng của bún thang được ninh từ xương gà và tôm khô, có vị ngọt thanh và thơm ngon. Đây là món ăn cầu kỳ, tinh tế và đòi hỏi nhiều công phu chế biến."
$ws.Range("A88").Value = "Miêu tả hương vị và quy trình chế biến của Chả Cá Thăng Long"
$ws.Range("B88").Value = "Chả Cá Thăng Long là món ăn đặc trưng của Hà Nội. Chả cá được ướp vàng ươm, nướng trên than hoa thơm lừng, ăn kèm bún, rau sống và mắm tôm. Miếng cá khi chín vàng hơi săn chắc, cháy cạnh giòn bên ngoài nhưng mềm ở giữa. Gia vị tẩm ướp vừa miệng và thơm phức. Khi ăn hòa quyện với mắm tôm tạo nên hương vị khó quên."
$ws.Range("A89").Value = "Hãy mô tả đặc điểm bánh cuốn tại Quán Bánh Cuốn Bà Xuân"
$ws.Range("B89").Value = "Bánh cuốn Bà Xuân có lớp bánh mỏng tang, trắng mịn, được tráng ngay tại chỗ và luôn nóng hổi. Nhân gồm thịt băm và mộc nhĩ thơm ngậy. Hành phi vàng giòn được rắc lên trên. Nước chấm pha chua ngọt hài hòa, ăn kèm vài lát chả quế và chút tinh dầu cà cuống, tạo nên hương vị tròn trịa, tinh tế."
$ws.Range("A90").Value = "Nêu đầy đủ mô tả về hương vị và trải nghiệm tại Quán Xôi Mây Vọng Đức"
$ws.Range("B90").Value = "Xôi Mây nổi tiếng gần 30 năm ở Hà Nội, có hương vị thơm ngon đặc trưng. Món xôi xéo có sắc vàng tươi, đậu xanh chín mềm, hành phi vàng sậm giòn béo nằm trong lớp lá chuối xanh biếc. Vào buổi sáng, nhiều người xếp hàng dài để mua xôi nóng hổi. Món ăn là sự kết hợp tinh tế giữa màu sắc và hương vị, được đánh giá là ngon nhất “vịnh Bắc Bộ”."
$ws.Range("A91").Value = "Miêu tả đầy đủ món ốc nấu chuối đậu tại Quán Ốc Dì Tú."
$ws.Range("B91").Value = "Tại Quán Ốc Dì Tú, món ốc nấu chuối đậu có ốc tươi rói, giòn sần sật, kết hợp với chuối đậu bùi bùi và nước dùng đậm đà. Hương vị độc đáo khiến thực khách muốn thưởng thức mãi không thôi. Đây là món ăn nổi bật và được nhiều thực khách yêu thích."
$ws.Range("A92").Value = "Bún ốc Hà Nội được làm từ những nguyên liệu chính nào?"
$ws.Range("B92").Value = "Bún ốc Hà Nội được làm từ ốc nhồi hoặc ốc vặn, bún tươi, cà chua chín mọng, rau thơm và gia vị dân dã. Nước dùng được ninh từ xương và nước luộc ốc, thêm giấm bỗng tạo vị chua nhẹ đặc trưng."
$ws.Range("A93").Value = "Địa chỉ bán phở Gia Truyền ở Hà Nội là đâu?"
$ws.Range("B93").Value = "Phở Gia Truyền có địa chỉ: 49 Bát Đàn, Quận Hoàn Kiếm, Hà Nội."
$ws.Range("A94").Value = "Miến lươn Hà Nội có những thành phần nào?"
$ws.Range("B94").Value = "Miến lươn được làm từ miến, lươn, hành lá, tía tô, rau răm, lạc rang, cùng nước dùng có vị ngọt, vị thơm ngậy của lươn, vị cay của ớt và hương thơm từ rau thơm."
$ws.Range("A95").Value = "Gà luộc Hà Nội khi ăn thường chấm với gì?"
$ws.Range("B95").Value = "Gà luộc Hà Nội thường được chấm với muối tiêu chanh hoặc mắm gừng."
$ws.Range("B96").Value = "Địa chỉ bán Kem Tràng Tiền: 35 Tràng Tiền, Quận Hoàn Kiếm, Hà Nội."
$ws.Range("A96").Value = "Địa chỉ bán Kem Tràng Tiền chỗ nào?"
$ws.Range("A97").Value = "Hãy mô tả hương vị và thành phần của món Phở Cuốn Hà Nội."
$ws.Range("B97").Value = "Phở cuốn gồm những lát bánh phở mềm cuộn thịt bò xào, rau sống, giá đỗ và dưa góp. Khi ăn được chấm với nước mắm pha tỏi ớt chua ngọt. Món ăn thanh mát, nhẹ nhàng, không cầu kỳ nhưng dễ gây nghiện."
$ws.Range("A98").Value = "Mô tả món Bún Mọc Dọc Mùng bao gồm những gì và hương vị ra sao?"
$ws.Range("B98").Value = "Bún mọc dọc mùng gồm viên mọc tròn thơm mùi mộc nhĩ, khoanh chân giò mềm rục, dọc mùng giòn sần sật và nước dùng trong, ngọt thanh từ xương ống ninh kỹ. Khi ăn thêm ớt cay, hành phi sẽ tạo trải nghiệm vị giác đậm chất Hà Thành, giản dị nhưng sâu sắc."
$ws.Range("A99").Value = "Hãy mô tả đầy đủ cách chế biến và hương vị của món Chả Cá Lã Vọng."
$ws.Range("B99").Value = "Chả cá Lã Vọng được làm từ cá lăng – loại cá ít xương, thịt ngọt và săn chắc. Cá được tẩm ướp gia vị đậm đà rồi chiên vàng trong chảo mỡ nóng. Khi ăn kết hợp với rau thì là, hành lá, bún tươi và mắm tôm pha. Hương vị béo, thơm, mằn mặn, dậy mùi đặc trưng, để lại dư vị khó quên."
$ws.Range("A100").Value = "Nêu mô tả đầy đủ về món sữa chua nếp cẩm Hà Nội."
$ws.Range("B100").Value = "Sữa chua nếp cẩm kết hợp vị chua dịu của sữa chua và vị ngọt bùi của nếp cẩm. Nếp có màu tím sẫm, ăn cùng sữa chua trắng mịn, thêm dừa khô và lạc rang giòn rụm. Hương vị gồm chút béo, chút thanh mát và dẻo bùi hòa quyện, là món ăn vặt quen thuộc của người Hà Nội."
$ws.Range("A101").Value = "bún riêu cua Hà Nội được chế biến như thế nào và có những thành phần gì?"
$ws.Range("B101").Value = "Bún riêu cua được làm từ cua đồng tươi giã nhuyễn, lọc lấy nước rồi đun nhỏ lửa để riêu cua kết thành từng mảng mềm. Nước dùng có cà chua tạo màu, thêm mắm tôm, đậu phụ rán vàng giòn, giò lụa, hành lá và rau sống. Hương vị đậm đà, có giấm bỗng chua thanh và ớt tươi tạo phong vị “rất Hà Nội”."

# Re-apply the standard wrap-text style to the newly populated rows (style #3)
$ws.Range("A69:B101").WrapText = $true

# WrapText triggers Excel's row auto-fit; pin every row back to the
# workbook's standard 49.95pt custom height (matches the target file).
$ws.Range("A69:B101").RowHeight = 49.95

# Restore the selection shown in the target workbook
$ws.Range("A55").Select() | Out-Null

